$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1168985
$ws.Range("C4").Value = 8211
$ws.Range("E4").Value = 927014
$ws.Range("F4").Value = 16478
$ws.Range("G4").Value = 510
$ws.Range("H4").Value = 67954

# Row 5 - España
$ws.Range("C5").Value = 1533

# Row 6 - Italia
$ws.Range("B6").Value = 210717
$ws.Range("C6").Value = 1389
$ws.Range("D6").Value = 81654
$ws.Range("E6").Value = 100179
$ws.Range("F6").Value = 1501
$ws.Range("G6").Value = 174
$ws.Range("H6").Value = 28884

# Row 9 - Alemania
$ws.Range("B9").Value = 165183
$ws.Range("C9").Value = 216
$ws.Range("E9").Value = 27771

# Row 11 - Turquia
$ws.Range("B11").Value = 126045
$ws.Range("C11").Value = 1670
$ws.Range("D11").Value = 63151
$ws.Range("E11").Value = 59497
$ws.Range("F11").Value = 1424
$ws.Range("G11").Value = 61
$ws.Range("H11").Value = 3397

# Row 21 - Ecuador
$ws.Range("B21").Value = 29538
$ws.Range("C21").Value = 2074
$ws.Range("D21").Value = 3300
$ws.Range("E21").Value = 24674
$ws.Range("G21").Value = 193
$ws.Range("H21").Value = 1564

# Row 56 - Marruecos
$ws.Range("B56").Value = 4903
$ws.Range("C56").Value = 174
$ws.Range("D56").Value = 1438
$ws.Range("E56").Value = 3291
